# The target edit lives entirely in word/styles.xml's <w:docDefaults> block
# (docDefaults isn't exposed as a discrete object in the Word OM), so we
# round-trip the whole package through Content.WordOpenXML, patch the
# docDefaults fragment textually, and write the package XML back.

$d = $word.ActiveDocument
$xml = $d.Content.WordOpenXML

$newDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

$start = $xml.IndexOf("<w:docDefaults>")
$end = $xml.IndexOf("</w:docDefaults>") + "</w:docDefaults>".Length
$oldDocDefaults = $xml.Substring($start, $end - $start)

$xml = $xml.Replace($oldDocDefaults, $newDocDefaults)

$d.Content.WordOpenXML = $xml
